$wb = $excel.ActiveWorkbook

function Set-ColumnValues($ws, $rangeAddr, $values) {
    $n = $values.Count
    $arr = New-Object 'object[,]' $n,1
    for ($i = 0; $i -lt $n; $i++) { $arr[$i,0] = $values[$i] }
    $ws.Range($rangeAddr).Value = $arr
}

function Set-RowValues($ws, $rangeAddr, $values) {
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range($rangeAddr).Value = $arr
}

# --- Sheet "TestCasesFlag": selection changes from A2:C6 to B2. ---
$flags = $wb.Worksheets.Item("TestCasesFlag")
$flags.Range("B2").Select()

# --- Sheet "Data": swap the placeholder columns (Value4..Value8 header /
# test2..test7 values) for the real columns used by the updated data
# provider function (password pair, name fields, address fields). ---
$data = $wb.Worksheets.Item("Data")

Set-ColumnValues $data "C2:C3" @("Password01", "Password02")
Set-RowValues    $data "D1:H1" @("FirstName", "LastName", "MidName ", "Add", "Add2")
Set-ColumnValues $data "D2:D3" @("AML admin", "AML controller")
Set-ColumnValues $data "E2:E3" @("UB", "PRABIN")
Set-ColumnValues $data "F2:F3" @("PRAKASH", "KAUR")
Set-ColumnValues $data "G2:G3" @("OSLO", "EDINBURG")
Set-ColumnValues $data "H2:H3" @("Done", "London")

# Active-cell selection on "Data" moves from H14 to C4 (Data stays the
# selected/active tab, as in the original workbook).
$data.Range("C4").Select()
